$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Pre-format the price cells whose new values look numeric (e.g. "1.000", "0.3764")
# as Text so Excel keeps the exact original string instead of coercing it to a number.
# (Done per-cell since a multi-area union Range only applies NumberFormat to its first area.)
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = '27.945.12'
$ws.Range("E2").Value = '  -1.45%  '

$ws.Range("D3").Value = '1.779.06'
$ws.Range("E3").Value = '  -1.34%  '

$ws.Range("E4").Value = '  +0.06%  '

$ws.Range("D5").Value = '315.92'
$ws.Range("E5").Value = '  +0.05%  '

$ws.Range("E6").Value = '  +0.10%  '

$ws.Range("D7").Value = '0.5387'
$ws.Range("E7").Value = '  -2.19%  '

$ws.Range("D8").Value = '0.3764'

$ws.Range("D9").Value = '0.07431'
$ws.Range("E9").Value = '  -2.16%  '

$ws.Range("D10").Value = '41.64'
$ws.Range("E10").Value = '  -2.35%  '

$ws.Range("E11").Value = '  -2.50%  '

$ws.Range("D12").Value = '1.000'
$ws.Range("E12").Value = '  +0.05%  '

$ws.Range("D13").Value = '20.42'
$ws.Range("E13").Value = '  -3.57%  '

$ws.Range("D14").Value = '6.069'
$ws.Range("E14").Value = '  -1.79%  '

$ws.Range("B15").Value = 'Chainlink'
$ws.Range("C15").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D15").Value = '7.199'
$ws.Range("E15").Value = '  -2.09%  '

$ws.Range("B16").Value = 'WrappedEther'
$ws.Range("C16").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D16").Value = '1.771.19'
$ws.Range("E16").Value = '  -1.72%  '

$ws.Range("D17").Value = '88.09'
$ws.Range("E17").Value = '  -4.44%  '

$ws.Range("D18").Value = '0.00001051'
$ws.Range("E18").Value = '  -1.49%  '

$ws.Range("D19").Value = '0.06437'
$ws.Range("E19").Value = '  +0.00%  '

$ws.Range("E20").Value = '  +0.07%  '

$ws.Range("D21").Value = '17.20'
$ws.Range("E21").Value = '  -0.80%  '

$ws.Range("E22").Value = '  -2.02%  '

$ws.Range("D23").Value = '27.979.26'
$ws.Range("E23").Value = '  -1.35%  '

$ws.Range("D24").Value = '11.14'
$ws.Range("E24").Value = '  -2.82%  '

$ws.Range("D25").Value = '2.088'
$ws.Range("E25").Value = '  -2.11%  '

$ws.Range("D26").Value = '156.03'
$ws.Range("E26").Value = '  -1.24%  '

$ws.Range("E27").Value = '  -2.33%  '

$ws.Range("D28").Value = '1.980.73'
$ws.Range("E28").Value = '  -1.44%  '

$ws.Range("D29").Value = '2.278'
$ws.Range("E29").Value = '  -4.62%  '

$ws.Range("D30").Value = '119.79'
$ws.Range("E30").Value = '  -3.25%  '

$ws.Range("D31").Value = '1.106'
$ws.Range("E31").Value = '  -1.54%  '

$ws.Range("D32").Value = '0.1046'
$ws.Range("E32").Value = '  +2.61%  '

$ws.Range("E33").Value = '  -0.77%  '

$ws.Range("D34").Value = '5.502'
$ws.Range("E34").Value = '  -4.19%  '

$ws.Range("E35").Value = '  -3.41%  '

$ws.Range("D36").Value = '0.06408'
$ws.Range("E36").Value = '  +1.11%  '

$ws.Range("D37").Value = '0.02264'
$ws.Range("E37").Value = '  -2.41%  '

$ws.Range("D38").Value = '4.960'

$ws.Range("D39").Value = '8.388'
$ws.Range("E39").Value = '  -5.33%  '

$ws.Range("D40").Value = '0.6125'
$ws.Range("E40").Value = '  -4.45%  '

$ws.Range("D41").Value = '11.05'
$ws.Range("E41").Value = '  -4.96%  '

$ws.Range("E42").Value = '  +3.36%  '

$ws.Range("D43").Value = '1.174'
$ws.Range("E43").Value = '  +1.42%  '

$ws.Range("D44").Value = '0.9996'
$ws.Range("E44").Value = '  +0.07%  '

$ws.Range("D45").Value = '13.28'
$ws.Range("E45").Value = '  -1.67%  '

$ws.Range("D46").Value = '3.660'
$ws.Range("E46").Value = '  -0.67%  '

$ws.Range("D47").Value = '0.5729'

$ws.Range("D48").Value = '126.32'
$ws.Range("E48").Value = '  +1.64%  '

$ws.Range("B49").Value = 'NEARProtocol'
$ws.Range("C49").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D49").Value = '1.919'
$ws.Range("E49").Value = '  -3.41%  '

$ws.Range("B50").Value = 'EOS'
$ws.Range("C50").Value = 'https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos'
$ws.Range("D50").Value = '1.178'
$ws.Range("E50").Value = '  +2.60%  '

$ws.Range("D51").Value = '0.06782'
$ws.Range("E51").Value = '  -1.79%  '
